$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet2" to "Sheet1" (also updates solver_* defined names
# that reference the sheet by name).
$ws.Name = "Sheet1"

# Row 4: B4 keeps its (empty) formatting but loses its border/alignment attributes;
# C4 and D4 are fully cleared out (no formatting left at all).
$ws.Range("C4:D4").Clear()

$wb.Save()
